$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Helper: range of a paragraph's content, excluding the trailing paragraph
# mark, so InsertXML only rewrites the runs and leaves the pilcrow's own
# pPr/rsid attributes on the <w:p> untouched.
function Get-ContentRange($para) {
    $r = $para.Range
    return $d.Range($r.Start, $r.End - 1)
}

# ---------------------------------------------------------------------
# Change 1: "(Dr. David Walker)" paragraph -> split "(" / "Dr." / " David
# Walker)" with spellcheck proofErr bookends around "Dr."
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$target = Get-ContentRange($p2)
$xml = "<w:p $wNs>" +
  "<w:r><w:rPr><w:b/><w:bCs/><w:u w:val=`"single`"/></w:rPr><w:t>(</w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:b/><w:bCs/><w:u w:val=`"single`"/></w:rPr><w:t>Dr.</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:b/><w:bCs/><w:u w:val=`"single`"/></w:rPr><w:t xml:space=`"preserve`"> David Walker)</w:t></w:r>" +
  "</w:p>"
$target.InsertXML($xml)

# ---------------------------------------------------------------------
# Change 2: "• Types of visualisations needed, graphs:" -> split into
# "• Types of visualisations " / "needed" / " graphs:" (same bold run
# formatting throughout).
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$target = Get-ContentRange($p5)
$xml = "<w:p $wNs>" +
  "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=`"preserve`">$([char]0x2022) Types of visualisations </w:t></w:r>" +
  "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>needed</w:t></w:r>" +
  "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=`"preserve`"> graphs:</w:t></w:r>" +
  "</w:p>"
$target.InsertXML($xml)

# ---------------------------------------------------------------------
# Change 3: "3D visualisations? ... multi-dimensional scaling" -> split
# so the trailing "multi-dimensional scaling" is highlighted yellow; the
# paragraph mark run properties also pick up the highlight (selection ran
# through end of paragraph).
# ---------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$target = Get-ContentRange($p6)
$xml = "<w:p $wNs>" +
  "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr>" +
  "<w:rPr><w:b/><w:bCs/><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:t xml:space=`"preserve`">3D visualisations? D3 does not out of the box have this, maybe some other solutions; </w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>multi-dimensional scaling</w:t></w:r>" +
  "</w:p>"
$target.InsertXML($xml)

# ---------------------------------------------------------------------
# Change 4: "make it modular" run loses its leading <w:tab/> and instead
# gets yellow highlighting.
# ---------------------------------------------------------------------
$p22 = $d.Paragraphs.Item(22)
$target = Get-ContentRange($p22)
$xml = "<w:p $wNs>" +
  "<w:r><w:t xml:space=`"preserve`">$([char]0x2022) </w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>make it modular</w:t></w:r>" +
  "</w:p>"
$target.InsertXML($xml)

# ---------------------------------------------------------------------
# Change 5: "• visualise pareto fronts" -> "• " + highlighted "visualise
# pareto fronts" + " (given the solutions from the optimiser)"
# ---------------------------------------------------------------------
$p23 = $d.Paragraphs.Item(23)
$target = Get-ContentRange($p23)
$xml = "<w:p $wNs>" +
  "<w:r><w:t xml:space=`"preserve`">$([char]0x2022) </w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>visualise pareto fronts</w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`"> (given the solutions from the optimiser)</w:t></w:r>" +
  "</w:p>"
$target.InsertXML($xml)

Write-Host "done"
